# slides14w: toggle the "hidden" flag on the Hat-Check / Chinese-Banquet
# slide run (slides 49-55). Slides 49-51 become hidden in the slide show;
# slides 52-55 (previously hidden) become shown again.

$p = $ppt.ActivePresentation

$toHide = @(49, 50, 51)
foreach ($idx in $toHide) {
    $slide = $p.Slides.Item($idx)
    $slide.SlideShowTransition.Hidden = $true
}

$toShow = @(52, 53, 54, 55)
foreach ($idx in $toShow) {
    $slide = $p.Slides.Item($idx)
    $slide.SlideShowTransition.Hidden = $false
}
